$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column C for rows 2-359
# from 2023-09-19 (serial 45188) to 2023-09-20 (serial 45189).
$ws.Range("C2:C359").Value = 45189
